# Updated cryptos list on Fri Sep 22 11:38:41 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows on the sheet. Values are written as literal text (a leading
# "'" forces Excel to keep decimal-looking numbers like "19.70" or "0.100"
# as text instead of normalizing them into numbers), matching how the
# source data is stored in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.669.45'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.598.11'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''211.21'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").Value = '''19.70'
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").Value = '''0.0838'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '1.821.38'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '1.587.33'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '''0.522'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '''64.82'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").Value = '26.649.65'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '''210.04'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '''2.30'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").Value = '''146.23'
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").Value = '''0.115'
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("D29").Value = '''15.31'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = '''0.0504'
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").Value = '''3.24'
$ws.Range("E32").Value = '  -0.69%  '
$ws.Range("D33").Value = '''0.670'
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").Value = '1.296.16'
$ws.Range("E35").Value = '  -1.18%  '
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("D37").Value = '''1.49'
$ws.Range("E37").Value = '  -1.63%  '
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("E42").Value = '  +1.15%  '
$ws.Range("D43").Value = '''0.787'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '''63.84'
$ws.Range("E44").Value = '  +1.87%  '
$ws.Range("D45").Value = '1.734.32'
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").Value = '''0.898'
$ws.Range("E46").Value = '  +10.98%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("D49").Value = '''0.100'
$ws.Range("E49").Value = '  +2.67%  '
$ws.Range("E50").Value = '  -1.15%  '
$ws.Range("D51").Value = '''7.52'
$ws.Range("E51").Value = '  +1.44%  '

# Reset style on cells that required a text-forcing apostrophe prefix,
# so Excel does not leave a stray quotePrefix style on them.
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
